$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.353225806451613
$ws.Range("C2").Value = 0.392338177014531
$ws.Range("D2").Value = 0.762430939226519
$ws.Range("E2").Value = 0.349378881987578
$ws.Range("F2").Value = 0.27884979190314

$ws.Range("B3").Value = 0.8
$ws.Range("C3").Value = 0.755614266842801
$ws.Range("D3").Value = 0.906077348066298
$ws.Range("E3").Value = 0.633540372670807
$ws.Range("F3").Value = 0.45251608021188

$ws.Range("B4").Value = 0.485483870967742
$ws.Range("C4").Value = 0.589167767503302
$ws.Range("D4").Value = 0.920810313075506
$ws.Range("E4").Value = 0.273291925465839
$ws.Range("F4").Value = 0.291335603480893
